$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("E2").Value = "'-1.08%"
$ws.Range("D3").Value = "'35.80"
$ws.Range("E3").Value = "'-0.36%"
$ws.Range("D4").Value = "'5.047"
$ws.Range("E4").Value = "'-0.35%"
$ws.Range("D5").Value = "'0.08044"
$ws.Range("E5").Value = "'-1.01%"
$ws.Range("D6").Value = "'1.863"
$ws.Range("E6").Value = "'-3.64%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.110"
$ws.Range("E7").Value = "'-1.15%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.771"
$ws.Range("E8").Value = "'-0.74%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9247"
$ws.Range("E9").Value = "'-1.56%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1281"
$ws.Range("E10").Value = "'-5.73%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1888"
$ws.Range("E11").Value = "'-1.07%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09030"
$ws.Range("E12").Value = "'-2.43%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03423"
$ws.Range("E13").Value = "'-2.42%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09864"
$ws.Range("E14").Value = "'-0.12%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001426"
$ws.Range("E15").Value = "'-1.72%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006243"
$ws.Range("E16").Value = "'7.69%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.866"
$ws.Range("E17").Value = "'6.74%"
$ws.Range("D18").Value = "'3.309"
$ws.Range("E18").Value = "'11.42%"
$ws.Range("D19").Value = "'0.3408"
$ws.Range("E19").Value = "'-1.26%"
$ws.Range("D20").Value = "'0.1333"
$ws.Range("E20").Value = "'-0.94%"
$ws.Range("D21").Value = "'4.812"
$ws.Range("E21").Value = "'-7.27%"
$ws.Range("E22").Value = "'-7.67%"
$ws.Range("E23").Value = "'-0.62%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-0.57%"
$ws.Range("D25").Value = "'0.004863"
$ws.Range("E25").Value = "'1.93%"
$ws.Range("D27").Value = "'0.0001300"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("E28").Value = "'42.22%"
$ws.Range("D39").Value = "'0.01961"
$ws.Range("E39").Value = "'-2.52%"
$ws.Range("D40").Value = "'0.05136"
$ws.Range("E40").Value = "'1.44%"
$ws.Range("D41").Value = "'0.007511"
$ws.Range("E41").Value = "'-1.37%"
$ws.Range("D42").Value = "'0.01018"
$ws.Range("E42").Value = "'-9.44%"
$ws.Range("E43").Value = "'-1.91%"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'0.54%"
$ws.Range("D45").Value = "'0.009893"
$ws.Range("E45").Value = "'-12.42%"
$ws.Range("D46").Value = "'0.00006103"
$ws.Range("E46").Value = "'-3.63%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.25%"
$ws.Range("D48").Value = "'64.85"
$ws.Range("E48").Value = "'-0.17%"
$ws.Range("D49").Value = "'0.001247"
$ws.Range("E49").Value = "'4.76%"
$ws.Range("D50").Value = "'0.00002094"
$ws.Range("E50").Value = "'-0.25%"
$ws.Range("D51").Value = "'0.0001994"
$ws.Range("E51").Value = "'-0.25%"
